$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
$ws.Range("C1").Value = "dist"
$ws.Range("D1").Value = "total"

# --- Data columns B (index) and C (dist) ---
$idx = @(0,1,2,3,4,5,6,7,8,9,10,11,$null,12,13,14,15,16,17,18)
$dist = @(0,1.36,2.41,2.73,1.2,1.4,1.6,1.84,3.76,0.54,2.9,1.84,0.16,1.6,1.4,1.2,2.65,2.64,1.04,7.01)

for ($i = 0; $i -lt 20; $i++) {
    $r = 2 + $i
    if ($r -eq 14) {
        $ws.Cells.Item($r, 2).Value = "11'"
    } else {
        $ws.Cells.Item($r, 2).Value = $idx[$i]
    }
    $ws.Cells.Item($r, 3).Value = $dist[$i]
}

# --- Column D formulas (running total) ---
$ws.Range("D2").Formula = "=C2"
for ($r = 3; $r -le 21; $r++) {
    $ws.Range("D$r").Formula = "=D$($r-1)+C$r"
}

# --- "ok" markers in column A ---
$ws.Range("A13").Value = "ok"
$ws.Range("A16").Value = "ok"
$ws.Range("A17").Value = "ok"

# --- Highlight rows 13 and 17 (A:D) with yellow fill ---
$ws.Range("A13:D13").Interior.Color = 65535
$ws.Range("A17:D17").Interior.Color = 65535

# --- Right-align the "11'" label in B14 ---
$ws.Range("B14").HorizontalAlignment = -4152

# --- Selection matching the saved view state ---
$ws.Range("D12:D21").Select()
